# Canary_UDP_Mapping.xlsx - "adapted canary to include ngchp changes"
#
# Logical changes to the PacketFormat sheet's NGCHP Gensets block:
#   1. Remove the "CHP NOx" row entirely.
#   2. Rename NGCHP "Fuel comsumption " row label to "Fuel usage ".
#   3. Add two new NGCHP rows at the end of the NGCHP block:
#        "Boiler Nm^3/hr"   (Units - Notes: Nm^3/hr)
#        "Boiler lbm/hr CO2" (Units - Notes: lbm/hr)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the "CHP NOx" row (row 18) -------------------------------
$ws.Rows.Item(18).Delete()

# --- 2. Rename "Fuel comsumption " (NGCHP, now still row 17) ------------
$ws.Range("B17").Value = "Fuel usage "

# --- 3. Insert two new rows right after "CHP heat contribution" ---------
#        (that row is now row 20; Motor data / F1 Motor Status is row 21)
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

# Fill in the new row 21: Boiler Nm^3/hr
$ws.Range("A21").Value = "NGCHP Gensets data"
$ws.Range("B21").Value = "Boiler Nm^3/hr"
$ws.Range("C21").Value = "int16"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("F21").Formula = "=RIGHT(C21,2)/8*D21*E21*`$C`$1"
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = "Nm^3/hr"

# Fill in the new row 22: Boiler lbm/hr CO2
$ws.Range("A22").Value = "NGCHP Gensets data"
$ws.Range("B22").Value = "Boiler lbm/hr CO2"
$ws.Range("C22").Value = "int16"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("F22").Formula = "=RIGHT(C22,2)/8*D22*E22*`$C`$1"
$ws.Range("G22").Value = 100
$ws.Range("H22").Value = "lbm/hr"

# --- 4. Re-assert the pre-existing (non-shared) formula on the "Power
#        Factor" row (originally row 24, now row 25) so its serialization
#        keeps standing out from the shared group exactly as it did before
#        this edit.
$ws.Range("F25").Formula = "=RIGHT(C25,2)/8*D25*E25*`$C`$1"

# --- 5. Restore the selection shown in the saved file --------------------
$ws.Range("C34").Select()
